{"js": "// Insert a new \"CONDITION / ICD CODE / TYPE / DATE OF ONSET\" detail block\n// right after the \"CLAIMED CONDITION\" heading, matching the canonical edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"CLAIMED CONDITION\" heading paragraph.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"CLAIMED CONDITION\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find \"CLAIMED CONDITION\" paragraph.');\n}\n\n// The rows to insert: [headingText, styleName, bodyText].\nconst rows = [\n  [\"CONDITION\", \"Heading 3\", \"Joint instability\"],\n  [\"ICD CODE\", \"Heading 3\", \"ICD-2017\"],\n  [\"TYPE\", \"Heading 3\", \"Accumulated over time (wear and tear)\"],\n  [\"DATE OF ONSET\", \"Heading 3\", \"2009-12-01\"],\n];\n\n// Insert in order, always directly after the growing \"target\" anchor so the\n// final order in the document matches the order of `rows`.\nlet anchor = target;\nfor (const [headingText, styleName, bodyText] of rows) {\n  const headingPara = anchor.insertParagraph(headingText, Word.InsertLocation.after);\n  headingPara.style = styleName;\n\n  const bodyPara = headingPara.insertParagraph(bodyText, Word.InsertLocation.after);\n  // New paragraphs inherit the preceding paragraph's style; reset this one\n  // back to the document default so it reads as plain body text.\n  bodyPara.style = \"Normal\";\n\n  anchor = bodyPara;\n}\n\nawait context.sync();\n", "ps1": "# Insert a \"CONDITION / ICD CODE / TYPE / DATE OF ONSET\" detail block right\n# after the \"CLAIMED CONDITION\" heading, matching the canonical edit.\n\n$d = $word.ActiveDocument\n\n# Locate the \"CLAIMED CONDITION\" heading paragraph.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`n\", \"`a\") -eq \"CLAIMED CONDITION\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find 'CLAIMED CONDITION' paragraph.\"\n}\n\n$rows = @(\n    @{ Heading = \"CONDITION\"; Body = \"Joint instability\" },\n    @{ Heading = \"ICD CODE\"; Body = \"ICD-2017\" },\n    @{ Heading = \"TYPE\"; Body = \"Accumulated over time (wear and tear)\" },\n    @{ Heading = \"DATE OF ONSET\"; Body = \"2009-12-01\" }\n)\n\n$anchorPara = $target\nforeach ($row in $rows) {\n    # Add an empty paragraph right after the anchor, then fill it in.\n    $anchorPara.Range.InsertParagraphAfter()\n    $headingPara = $anchorPara.Next()\n    $headingPara.Range.InsertBefore($row.Heading)\n    $headingPara.Range.Style = \"Heading 3\"\n\n    $headingPara.Range.InsertParagraphAfter()\n    $bodyPara = $headingPara.Next()\n    $bodyPara.Range.InsertBefore($row.Body)\n    $bodyPara.Range.Style = \"Normal\"\n\n    $anchorPara = $bodyPara\n}\n"}
